$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume/1h change (E) figures per the latest data refresh.
# Cells whose new text would otherwise be auto-parsed as a number by Excel are
# temporarily switched to text format ("@") so the value is stored as a string,
# matching the original inline-string cell contents, then restored to General format.

$ws.Range('D2').Value = '27.167.40'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '1.847.24'
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4632'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3703'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07376'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8836'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07938'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.93'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('D13').Value = '1.913.44'
$ws.Range('E13').Value = '  +5.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.377'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.590'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008946'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.08%  '
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.16%  '
$ws.Range('D21').Value = '27.194.92'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.141'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').Value = '2.102.81'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.866'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.139'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08890'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.972'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7416'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.23%  '
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.549'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.26%  '
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05276'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01952'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.966'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.096'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5173'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1638'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.282'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4860'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.27'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.81'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06232'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.34%  '
